# Updated cryptos list on Mon May  8 10:48:44 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with a new snapshot of values. Every cell on this sheet is stored
# as literal text (not a number), including price strings that use a dot
# as a thousands separator (e.g. "27.960.92") and percent strings padded
# with surrounding spaces (e.g. "  -3.26%  "). Plain decimal-looking price
# values (e.g. "318.10") would otherwise be auto-converted to a numeric
# type by Excel's normal smart-entry parsing, which would silently drop
# the original text formatting - so those are entered with a leading
# apostrophe to force text, exactly like typing them by hand in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.960.92'
$ws.Range("E2").Value = '  -3.26%  '

$ws.Range("D3").Value = '1.859.70'
$ws.Range("E3").Value = '  -2.52%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''318.10'
$ws.Range("E5").Value = '  -2.03%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("E7").Value = '  -4.88%  '

$ws.Range("D8").Value = '''0.3697'
$ws.Range("E8").Value = '  -3.11%  '

$ws.Range("D9").Value = '''0.07497'
$ws.Range("E9").Value = '  -2.98%  '

$ws.Range("D10").Value = '''0.9418'
$ws.Range("E10").Value = '  -3.85%  '

$ws.Range("D11").Value = '''21.30'
$ws.Range("E11").Value = '  -3.52%  '

$ws.Range("D12").Value = '1.840.12'
$ws.Range("E12").Value = '  -3.80%  '

$ws.Range("D13").Value = '''6.728'
$ws.Range("E13").Value = '  -3.08%  '

$ws.Range("D14").Value = '''5.437'
$ws.Range("E14").Value = '  -4.24%  '

$ws.Range("D15").Value = '''0.06851'
$ws.Range("E15").Value = '  -3.03%  '

$ws.Range("D16").Value = '''1.003'
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").Value = '''81.61'
$ws.Range("E17").Value = '  -2.49%  '

$ws.Range("D18").Value = '''0.000009025'
$ws.Range("E18").Value = '  -4.57%  '

$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("E20").Value = '  -3.94%  '

$ws.Range("D21").Value = '27.932.31'
$ws.Range("E21").Value = '  -3.35%  '

$ws.Range("D22").Value = '''5.115'
$ws.Range("E22").Value = '  -3.77%  '

$ws.Range("D23").Value = '''11.08'
$ws.Range("E23").Value = '  +1.32%  '

$ws.Range("D24").Value = '2.100.31'
$ws.Range("E24").Value = '  -2.04%  '

$ws.Range("E25").Value = '  -4.31%  '

$ws.Range("D26").Value = '''154.71'
$ws.Range("E26").Value = '  -2.46%  '

$ws.Range("D27").Value = '''18.38'
$ws.Range("E27").Value = '  -3.61%  '

$ws.Range("D28").Value = '''5.380'
$ws.Range("E28").Value = '  -5.01%  '

$ws.Range("D29").Value = '''113.68'
$ws.Range("E29").Value = '  -3.21%  '

$ws.Range("D30").Value = '''1.734'
$ws.Range("E30").Value = '  -7.23%  '

$ws.Range("D31").Value = '''0.08981'
$ws.Range("E31").Value = '  -3.30%  '

$ws.Range("D32").Value = '''0.8102'
$ws.Range("E32").Value = '  -5.96%  '

$ws.Range("D33").Value = '''4.818'
$ws.Range("E33").Value = '  -5.27%  '

$ws.Range("E34").Value = '  -5.86%  '

$ws.Range("D35").Value = '''2.938'
$ws.Range("E35").Value = '  -3.36%  '

$ws.Range("D36").Value = '''1.001'
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").Value = '''0.05482'
$ws.Range("E37").Value = '  -3.97%  '

$ws.Range("E38").Value = '  -3.82%  '

$ws.Range("D39").Value = '''0.01972'
$ws.Range("E39").Value = '  -3.47%  '

$ws.Range("D40").Value = '''2.919'
$ws.Range("E40").Value = '  +2.08%  '

$ws.Range("D41").Value = '''0.5252'

$ws.Range("D42").Value = '''7.012'
$ws.Range("E42").Value = '  -5.47%  '

$ws.Range("D43").Value = '''0.1685'
$ws.Range("E43").Value = '  -3.75%  '

$ws.Range("D44").Value = '''8.800'
$ws.Range("E44").Value = '  -5.51%  '

$ws.Range("D45").Value = '''0.06795'
$ws.Range("E45").Value = '  -1.53%  '

$ws.Range("E46").Value = '  -5.45%  '

$ws.Range("D47").Value = '''10.60'
$ws.Range("E47").Value = '  -5.44%  '

$ws.Range("D48").Value = '''106.17'
$ws.Range("E48").Value = '  -3.73%  '

$ws.Range("D49").Value = '''1.680'
$ws.Range("E49").Value = '  -5.48%  '

$ws.Range("D50").Value = '''1.908'
$ws.Range("E50").Value = '  -10.19%  '

$ws.Range("D51").Value = '''1.000'
$ws.Range("E51").Value = '  -0.13%  '
